$wb = $excel.ActiveWorkbook

# --- TestScriptMasterSheet: row5 "testM3"/"Yes" -> "testTm3"/"YES" ---
$wsMaster = $wb.Worksheets.Item("TestScriptMasterSheet")
$wsMaster.Range("C5").Value = "YES"
$wsMaster.Range("A5").Value = "testTm3"

# --- Login: row5 "testM3" -> "testTm3" ---
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("A5").Value = "testTm3"

# --- ValidatePerson: row2 "testM3" -> "testTm3" ---
$wsValidate = $wb.Worksheets.Item("ValidatePerson")
$wsValidate.Range("A2").Value = "testTm3"

# --- Update the on-screen selection on each sheet (tab order), ending on
#     ValidatePerson so it stays the active sheet/tab, matching the diff ---
$wsLogin.Range("A5").Select() | Out-Null

$wsScreening = $wb.Worksheets.Item("Screening")
$wsScreening.Range("A2").Select() | Out-Null

$wsFolio = $wb.Worksheets.Item("Folio")
$wsFolio.Range("A3").Select() | Out-Null

$wsServices = $wb.Worksheets.Item("Services")
$wsServices.Range("A3").Select() | Out-Null

$wsValidate.Range("A2").Select() | Out-Null
